$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column C first (NewBuildingName / building6 / building10),
# then column D (NewFloorname / floorNo12 / FloorNo5) so shared strings
# are appended in the same order as the target workbook.
$ws.Range("C1").Value = "NewBuildingName"
$ws.Range("C2").Value = "building6"
$ws.Range("C3").Value = "building10"

$ws.Range("D1").Value = "NewFloorname"
$ws.Range("D2").Value = "floorNo12"
$ws.Range("D3").Value = "FloorNo5"

# Give the new header cells (C1/D1) the same (yellow-fill) style as the
# existing header cells A1/B1.
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)

# New column widths for C and D.
$ws.Columns.Item(3).ColumnWidth = 16.33
$ws.Columns.Item(4).ColumnWidth = 17

# Selection moves to B1.
[void]$ws.Range("B1").Select()
